$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'"
$ws.Range("C2").Style = "Normal"
$ws.Range("F2").Value = 36.8
$ws.Range("G2").Value = 42.2
$ws.Range("H2").Value = 40.4
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 33
$ws.Range("O2").Value = 40
$ws.Range("P2").Value = 22
$ws.Range("Q2").Value = 39
$ws.Range("R2").Value = 50
$ws.Range("AJ2").Value = 0.1
$ws.Range("AK2").Value = 6
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 0
$ws.Range("A5").Value = "PHI"
$ws.Range("B5").Value = "Tyrese Maxey"
$ws.Range("C5").Value = "'"
$ws.Range("C5").Style = "Normal"
$ws.Range("E5").Value = "'"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = 37.8
$ws.Range("G5").Value = 34.7
$ws.Range("H5").Value = 35.9
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = "-"
$ws.Range("O5").Value = "-"
$ws.Range("P5").Value = "-"
$ws.Range("Q5").Value = "-"
$ws.Range("R5").Value = 32
$ws.Range("T5").Value = "NYK"
$ws.Range("V5").Value = 46
$ws.Range("W5").Value = "vs"
$ws.Range("X5").Value = 41
$ws.Range("AA5").Value = "@"
$ws.Range("AB5").Value = "MIL"
$ws.Range("AC5").Value = "vs"
$ws.Range("AD5").Value = "CHA"
$ws.Range("AF5").Value = "MIA"
$ws.Range("AG5").Value = "@"
$ws.Range("AH5").Value = "PHX"
$ws.Range("AJ5").Value = 0.4
$ws.Range("AK5").Value = "'"
$ws.Range("AK5").Style = "Normal"
$ws.Range("AL5").Value = "'"
$ws.Range("AL5").Style = "Normal"
$ws.Range("A7").Value = "BOS"
$ws.Range("B7").Value = "Jaylen Brown"
$ws.Range("D7").Value = "G"
$ws.Range("E7").Value = "O"
$ws.Range("F7").Value = 37
$ws.Range("H7").Value = 31.4
$ws.Range("I7").Value = 10
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 30
$ws.Range("O7").Value = 32
$ws.Range("P7").Value = 59
$ws.Range("Q7").Value = 24
$ws.Range("R7").Value = 40
$ws.Range("T7").Value = "UTA"
$ws.Range("U7").Value = "vs"
$ws.Range("V7").Value = 18
$ws.Range("AB7").Value = "PHX"
$ws.Range("AC7").Value = "@"
$ws.Range("AD7").Value = "WAS"
$ws.Range("AF7").Value = "DET"
$ws.Range("AG7").Value = "vs"
$ws.Range("AH7").Value = "MIL"
$ws.Range("AJ7").Value = -1.7
$ws.Range("AK7").Value = -0.7
$ws.Range("AL7").Value = 8
$ws.Range("A8").Value = "IND"
$ws.Range("B8").Value = "Pascal Siakam"
$ws.Range("D8").Value = "F"
$ws.Range("F8").Value = 28.4
$ws.Range("G8").Value = 33.1
$ws.Range("H8").Value = 34.2
$ws.Range("I8").Value = 11
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 2
$ws.Range("N8").Value = 32
$ws.Range("O8").Value = 30
$ws.Range("P8").Value = 34
$ws.Range("Q8").Value = 31
$ws.Range("R8").Value = 15
$ws.Range("S8").Value = "@"
$ws.Range("T8").Value = "OKC"
$ws.Range("U8").Value = "-"
$ws.Range("V8").Value = "-"
$ws.Range("W8").Value = "-"
$ws.Range("X8").Value = "-"
$ws.Range("Y8").Value = "-"
$ws.Range("Z8").Value = "-"
$ws.Range("AA8").Value = "vs"
$ws.Range("AC8").Value = "vs"
$ws.Range("AD8").Value = "BKN"
$ws.Range("AF8").Value = "CLE"
$ws.Range("AH8").Value = "DET"
$ws.Range("AI8").Value = "@"
$ws.Range("AJ8").Value = -1.4
$ws.Range("A9").Value = "LAC"
$ws.Range("B9").Value = "James Harden"
$ws.Range("D9").Value = "G"
$ws.Range("F9").Value = 27.8
$ws.Range("G9").Value = 32.8
$ws.Range("H9").Value = 31.6
$ws.Range("J9").Value = 3
$ws.Range("M9").Value = 3
$ws.Range("N9").Value = 16
$ws.Range("O9").Value = 33
$ws.Range("P9").Value = 37
$ws.Range("Q9").Value = 44
$ws.Range("R9").Value = 9
$ws.Range("T9").Value = "MIN"
$ws.Range("U9").Value = "@"
$ws.Range("V9").Value = 9
$ws.Range("W9").Value = "vs"
$ws.Range("X9").Value = 25
$ws.Range("Y9").Value = "@"
$ws.Range("Z9").Value = 15
$ws.Range("AA9").Value = "@"
$ws.Range("AB9").Value = "CHI"
$ws.Range("AD9").Value = "NOP"
$ws.Range("AF9").Value = "ATL"
$ws.Range("AH9").Value = "POR"
$ws.Range("AJ9").Value = 0.4
$ws.Range("A10").Value = "OKC"
$ws.Range("B10").Value = "Chet Holmgren"
$ws.Range("F10").Value = 22.4
$ws.Range("G10").Value = 32.3
$ws.Range("H10").Value = 30.6
$ws.Range("I10").Value = 11
$ws.Range("L10").Value = 3
$ws.Range("M10").Value = 4
$ws.Range("N10").Value = 9
$ws.Range("O10").Value = 8
$ws.Range("P10").Value = 47
$ws.Range("Q10").Value = 21
$ws.Range("R10").Value = 27
$ws.Range("T10").Value = "IND"
$ws.Range("U10").Value = "-"
$ws.Range("V10").Value = "-"
$ws.Range("AB10").Value = "DAL"
$ws.Range("AC10").Value = "@"
$ws.Range("AD10").Value = "MEM"
$ws.Range("AE10").Value = "vs"
$ws.Range("AF10").Value = "UTA"
$ws.Range("AH10").Value = "TOR"
$ws.Range("AJ10").Value = -0.2
$ws.Range("C14").Value = "Questionable"
$ws.Range("I18").Value = 7
$ws.Range("M18").Value = 3
$ws.Range("O18").Value = "-"
$ws.Range("P18").Value = 30
$ws.Range("Q18").Value = 43
$ws.Range("R18").Value = "-"
$ws.Range("A19").Value = "MEM"
$ws.Range("B19").Value = "Jaren Jackson Jr."
$ws.Range("D19").Value = "C"
$ws.Range("F19").Value = 32.6
$ws.Range("G19").Value = 31.1
$ws.Range("H19").Value = 30.5
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 2
$ws.Range("M19").Value = 2
$ws.Range("N19").Value = 18
$ws.Range("O19").Value = 33
$ws.Range("P19").Value = 43
$ws.Range("R19").Value = "-"
$ws.Range("S19").Value = "vs"
$ws.Range("T19").Value = "WAS"
$ws.Range("U19").Value = "@"
$ws.Range("V19").Value = 26
$ws.Range("W19").Value = "-"
$ws.Range("X19").Value = "-"
$ws.Range("AA19").Value = "vs"
$ws.Range("AB19").Value = "CHA"
$ws.Range("AD19").Value = "OKC"
$ws.Range("AE19").Value = "@"
$ws.Range("AF19").Value = "SAC"
$ws.Range("AH19").Value = "GSW"
$ws.Range("AI19").Value = "vs"
$ws.Range("AJ19").Value = -2.3
